# Adds the "Human Resources Attestation" scenario results (rows 45-90)
# to the SenaryoSonuclari (Scenario Results) worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$title = "The Admin Add Edit and Delete Attestations under Human Resources"

# Pass/Fail pattern for column B, rows 45-90 (row number => result)
$results = @{
    45 = "FAILED"; 46 = "PASSED"; 47 = "PASSED"; 48 = "PASSED"; 49 = "PASSED";
    50 = "FAILED"; 51 = "PASSED"; 52 = "FAILED"; 53 = "FAILED"; 54 = "FAILED";
    55 = "FAILED"; 56 = "FAILED"; 57 = "FAILED"; 58 = "FAILED"; 59 = "FAILED";
    60 = "FAILED"; 61 = "FAILED"; 62 = "FAILED"; 63 = "FAILED"; 64 = "FAILED";
    65 = "FAILED"; 66 = "FAILED"; 67 = "FAILED"; 68 = "FAILED"; 69 = "PASSED";
    70 = "PASSED"; 71 = "FAILED"; 72 = "FAILED"; 73 = "FAILED"; 74 = "FAILED";
    75 = "PASSED"; 76 = "PASSED"; 77 = "PASSED"; 78 = "FAILED"; 79 = "FAILED";
    80 = "FAILED"; 81 = "FAILED"; 82 = "FAILED"; 83 = "FAILED"; 84 = "FAILED";
    85 = "PASSED"; 86 = "PASSED"; 87 = "PASSED"; 88 = "PASSED"; 89 = "PASSED";
    90 = "PASSED"
}

for ($r = 45; $r -le 90; $r++) {
    $ws.Cells.Item($r, 1).Value = $title
    $ws.Cells.Item($r, 2).Value = $results[$r]
    $ws.Cells.Item($r, 3).Value = "chrome"
}
